# Weekly update: insert a new price record (row 240) for
# "Vega Central Mapocho de Santiago - Ciboulette" and push the
# existing rows (240-276) down by one (to 241-277).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 240; this shifts the previous
# rows 240..276 down to 241..277 and carries their formatting
# (including the date number format on column D) along with them.
$ws.Rows.Item(240).Insert()

# Populate the newly inserted row 240 with the new weekly record.
$ws.Range("A240").Value = 9
$ws.Range("B240").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C240").Value = "Metropolitana"
$ws.Range("D240").Value = 44491
$ws.Range("E240").Value = 13
$ws.Range("F240").Value = 100112039
$ws.Range("G240").Value = "Ciboulette"
$ws.Range("H240").Value = "Sin especificar"
$ws.Range("I240").Value = "Primera"
$ws.Range("J240").Value = 250
$ws.Range("K240").Value = 800
$ws.Range("L240").Value = 1000
$ws.Range("M240").Value = 900
$ws.Range("N240").Value = "`$/docena de atados"
$ws.Range("O240").Value = "Región Metropolitana"
$ws.Range("P240").Value = 300
$ws.Range("Q240").Value = 3
$ws.Range("R240").Value = "Hortaliza"
